$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.001.82"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "2.652.73"
$ws.Range("E3").Value = "  +1.63%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'579.58"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'145.18"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.23%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("D10").Value = "'0.107"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  +5.01%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "3.124.08"
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "'25.90"
$ws.Range("E14").Value = "  +11.29%  "
$ws.Range("D15").Value = "60.985.13"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "2.665.74"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "'11.69"
$ws.Range("E18").Value = "  +3.18%  "
$ws.Range("E19").Value = "  +2.41%  "
$ws.Range("D20").Value = "'351.21"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").Value = "'0.539"
$ws.Range("E22").Value = "  +2.54%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'63.96"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.163"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D28").Value = "'1.97"
$ws.Range("E28").Value = "  +5.66%  "
$ws.Range("E29").Value = "  +2.69%  "
$ws.Range("D30").Value = "'6.88"
$ws.Range("E30").Value = "  +7.78%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'165.42"
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").Value = "'19.97"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("E34").Value = "  +10.33%  "
$ws.Range("E35").Value = "  +5.78%  "
$ws.Range("E36").Value = "  +6.20%  "
$ws.Range("E37").Value = "  +3.46%  "
$ws.Range("D38").Value = "'333.94"
$ws.Range("E38").Value = "  +11.93%  "
$ws.Range("D39").Value = "'4.03"
$ws.Range("E39").Value = "  +4.17%  "
$ws.Range("D40").Value = "'38.70"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'5.23"
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'20.49"
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("D44").Value = "'134.61"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").Value = "'0.0564"
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("E47").Value = "  +3.17%  "
$ws.Range("D48").Value = "'0.616"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").Value = "'20.54"
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "2.098.57"
$ws.Range("E51").Value = "  +3.49%  "
